# Add a new "2022-Q3" quarter snapshot to the workbook.
#  1. Insert a new first-data-row on the "总计" (summary) sheet for 2022-Q3.
#  2. Insert a brand-new "2022-Q3" worksheet (before "2022-Q2") holding the
#     per-fund holding breakdown for the new quarter.
# All the older quarter sheets simply shift one position to the right and
# keep their own name/content unchanged.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) 总计 (summary) sheet: insert a new row 2 for 2022-Q3
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

# The freshly inserted row inherited some stray formatting from the
# insert operation - put every cell back the way the rest of the table
# looks: col A uses the bordered/centered "index" style, B:D are plain.
$summary.Cells.Item(3,1).Copy($summary.Cells.Item(2,1))
$summary.Cells.Item(2,2).Style = "Normal"
$summary.Cells.Item(2,3).Style = "Normal"
$summary.Cells.Item(2,4).Style = "Normal"

$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 8
$summary.Cells.Item(2,4).Value = 1.18

# Column A is just a running 0-based row index - renumber rows 2..9 now
# that an extra row was spliced in at the top.
for ($i = 0; $i -le 7; $i++) {
    $summary.Cells.Item(2 + $i, 1).Value = $i
}

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet with the per-fund breakdown.
#    Duplicate the "2022-Q2" sheet (same layout/styling) and drop it in
#    right before "2022-Q2", then overwrite its figures with the Q3 data.
# ---------------------------------------------------------------------
$prevQ = $wb.Worksheets.Item("2022-Q2")
$prevQ.Copy($prevQ)
$q3 = $wb.ActiveSheet
$q3.Name = "2022-Q3"

# Row 2 - 159869 华夏中证动漫游戏ETF (code/name unchanged)
$q3.Cells.Item(2,4).Value = "'6.35"
$q3.Cells.Item(2,5).Value = "'99.31"
$q3.Cells.Item(2,6).Value = "'9.38"
$q3.Cells.Item(2,7).Value = "'0.5956"
$q3.Cells.Item(2,8).Value = 3

# Row 3 - 516010 国泰中证动漫游戏ETF (code/name unchanged)
$q3.Cells.Item(3,4).Value = "'3.78"
$q3.Cells.Item(3,5).Value = "'97.86"
$q3.Cells.Item(3,6).Value = "'9.09"
$q3.Cells.Item(3,7).Value = "'0.3436"
$q3.Cells.Item(3,8).Value = 3

# Row 4 - 516770 华泰柏瑞中证动漫游戏ETF (code/name unchanged)
$q3.Cells.Item(4,4).Value = "'0.99"
$q3.Cells.Item(4,5).Value = "'96.39"
$q3.Cells.Item(4,6).Value = "'9.26"
$q3.Cells.Item(4,7).Value = "'0.0917"
$q3.Cells.Item(4,8).Value = 2

# Row 5 - 161030, name gains an "A" share-class suffix
$q3.Cells.Item(5,3).Value = "富国中证体育产业指数A"
$q3.Cells.Item(5,4).Value = "'1.59"
$q3.Cells.Item(5,5).Value = "'94.00"
$q3.Cells.Item(5,6).Value = "'4.41"
$q3.Cells.Item(5,7).Value = "'0.0701"
$q3.Cells.Item(5,8).Value = 8

# Row 6 - new fund: 517500 国泰中证沪港深动漫游戏ETF
$q3.Cells.Item(6,2).Value = "'517500"
$q3.Cells.Item(6,3).Value = "国泰中证沪港深动漫游戏ETF"
$q3.Cells.Item(6,4).Value = "'0.53"
$q3.Cells.Item(6,5).Value = "'92.78"
$q3.Cells.Item(6,6).Value = "'6.00"
$q3.Cells.Item(6,7).Value = "'0.0318"
$q3.Cells.Item(6,8).Value = 4

# Row 7 - new fund: 013278 富国中证体育产业指数C
$q3.Cells.Item(7,2).Value = "'013278"
$q3.Cells.Item(7,3).Value = "富国中证体育产业指数C"
$q3.Cells.Item(7,4).Value = "'0.42"
$q3.Cells.Item(7,5).Value = "'94.00"
$q3.Cells.Item(7,6).Value = "'4.41"
$q3.Cells.Item(7,7).Value = "'0.0185"
$q3.Cells.Item(7,8).Value = 8

# Row 8 - new fund: 159725 工银瑞信中证线上消费主题ETF
$q3.Cells.Item(8,2).Value = "'159725"
$q3.Cells.Item(8,3).Value = "工银瑞信中证线上消费主题ETF"
$q3.Cells.Item(8,4).Value = "'0.57"
$q3.Cells.Item(8,5).Value = "'98.42"
$q3.Cells.Item(8,6).Value = "'3.24"
$q3.Cells.Item(8,7).Value = "'0.0185"
$q3.Cells.Item(8,8).Value = 9

# Row 9 - new fund: 159728 南方国证在线消费ETF
$q3.Cells.Item(9,2).Value = "'159728"
$q3.Cells.Item(9,3).Value = "南方国证在线消费ETF"
$q3.Cells.Item(9,4).Value = "'0.33"
$q3.Cells.Item(9,5).Value = "'99.99"
$q3.Cells.Item(9,6).Value = "'3.33"
$q3.Cells.Item(9,7).Value = "'0.0110"
$q3.Cells.Item(9,8).Value = 7
